$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 5559.7
$ws.Cells.Item(64, 9).Value = 4599.6665
$ws.Cells.Item(64, 11).Value = 4599.6665
$ws.Cells.Item(64, 13).Value = -4351.6665
$ws.Cells.Item(67, 8).Value = 5559.7
$ws.Cells.Item(67, 9).Value = 4599.6665
$ws.Cells.Item(67, 11).Value = 4599.6665
$ws.Cells.Item(67, 13).Value = -3741.6665
$ws.Cells.Item(69, 8).Value = 7010.5713
$ws.Cells.Item(69, 10).Value = 7010.5713
$ws.Cells.Item(69, 12).Value = 21031.7139
$ws.Cells.Item(69, 14).Value = -22779.7139
$ws.Cells.Item(72, 8).Value = 7010.5713
$ws.Cells.Item(72, 10).Value = 7010.5713
$ws.Cells.Item(72, 12).Value = 63095.14169999999
$ws.Cells.Item(72, 14).Value = -71831.14169999999
$ws.Cells.Item(88, 8).Value = 772038.9
$ws.Cells.Item(88, 9).Value = 3213.3333
$ws.Cells.Item(88, 10).Value = 1002686.5
$ws.Cells.Item(88, 11).Value = 3213.3333
$ws.Cells.Item(88, 12).Value = 1002686.5
$ws.Cells.Item(88, 13).Value = -2807.3333
$ws.Cells.Item(88, 14).Value = -1003498.5
$ws.Cells.Item(91, 8).Value = 772038.9
$ws.Cells.Item(91, 9).Value = 3213.3333
$ws.Cells.Item(91, 10).Value = 1002686.5
$ws.Cells.Item(91, 11).Value = 3213.3333
$ws.Cells.Item(91, 12).Value = 1002686.5
$ws.Cells.Item(91, 13).Value = -1809.3333
$ws.Cells.Item(91, 14).Value = -1005494.5
$ws.Cells.Item(100, 8).Value = 2798.2778
$ws.Cells.Item(100, 10).Value = 3251
$ws.Cells.Item(100, 12).Value = 3251
$ws.Cells.Item(100, 14).Value = -4333
$ws.Cells.Item(112, 8).Value = 2328.2693
$ws.Cells.Item(112, 10).Value = 2360
$ws.Cells.Item(112, 12).Value = 7080
$ws.Cells.Item(112, 14).Value = -9296
$ws.Cells.Item(129, 8).Value = 2170.8333
$ws.Cells.Item(129, 9).Value = 1438.4286
$ws.Cells.Item(129, 10).Value = 3196.2
$ws.Cells.Item(129, 11).Value = 4315.2858
$ws.Cells.Item(129, 12).Value = 9588.599999999999
$ws.Cells.Item(129, 13).Value = 684.7142000000003
$ws.Cells.Item(129, 14).Value = -19588.6
$ws.Cells.Item(132, 8).Value = 5057.3447
$ws.Cells.Item(132, 9).Value = 2756.5
$ws.Cells.Item(132, 11).Value = 8269.5
$ws.Cells.Item(132, 13).Value = -5739.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3114.4375
$ws.Cells.Item(2, 9).Value = 1213.1818
$ws.Cells.Item(2, 11).Value = 1213.1818
$ws.Cells.Item(2, 13).Value = -1100.1818
$ws.Cells.Item(43, 8).Value = 32457.334
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 13).Value = $null
$ws.Cells.Item(61, 8).Value = 5963
$ws.Cells.Item(61, 9).Value = 5588.6
$ws.Cells.Item(61, 11).Value = 5588.6
$ws.Cells.Item(61, 13).Value = -5376.6
$ws.Cells.Item(97, 8).Value = 819.46155
$ws.Cells.Item(97, 10).Value = 1439.6
$ws.Cells.Item(97, 12).Value = 1439.6
$ws.Cells.Item(97, 14).Value = -2431.6
$ws.Cells.Item(116, 8).Value = 3114.4375
$ws.Cells.Item(116, 9).Value = 1213.1818
$ws.Cells.Item(116, 11).Value = 1213.1818
$ws.Cells.Item(116, 13).Value = 1080.8182
$ws.Cells.Item(122, 8).Value = 2941.138
$ws.Cells.Item(122, 9).Value = 3988.6
$ws.Cells.Item(122, 10).Value = 1818.8572
$ws.Cells.Item(122, 11).Value = 11965.8
$ws.Cells.Item(122, 12).Value = 5456.571599999999
$ws.Cells.Item(122, 13).Value = -9515.799999999999
$ws.Cells.Item(122, 14).Value = -10356.5716
$ws.Cells.Item(136, 8).Value = 5963
$ws.Cells.Item(136, 9).Value = 5588.6
$ws.Cells.Item(136, 11).Value = 16765.8
$ws.Cells.Item(136, 13).Value = -14215.8
$ws.Cells.Item(138, 8).Value = 99941.75
$ws.Cells.Item(138, 10).Value = 99941.75
$ws.Cells.Item(138, 12).Value = 99941.75
$ws.Cells.Item(138, 14).Value = -110221.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3114.4375
$ws.Cells.Item(3, 9).Value = 1213.1818
$ws.Cells.Item(3, 11).Value = 1213.1818
$ws.Cells.Item(3, 13).Value = -1099.1818
$ws.Cells.Item(62, 8).Value = 68450
$ws.Cells.Item(62, 10).Value = 60000
$ws.Cells.Item(62, 12).Value = 60000
$ws.Cells.Item(62, 14).Value = -61372
$ws.Cells.Item(65, 8).Value = 68450
$ws.Cells.Item(65, 10).Value = 60000
$ws.Cells.Item(65, 12).Value = 180000
$ws.Cells.Item(65, 14).Value = -186864
$ws.Cells.Item(75, 8).Value = 7491.3335
$ws.Cells.Item(75, 9).Value = 7491.3335
$ws.Cells.Item(75, 11).Value = 7491.3335
$ws.Cells.Item(75, 13).Value = -6555.3335
$ws.Cells.Item(78, 8).Value = 7491.3335
$ws.Cells.Item(78, 9).Value = 7491.3335
$ws.Cells.Item(78, 11).Value = 22474.0005
$ws.Cells.Item(78, 13).Value = -17794.0005
$ws.Cells.Item(94, 8).Value = 1160.2222
$ws.Cells.Item(94, 9).Value = 1192.75
$ws.Cells.Item(94, 11).Value = 1192.75
$ws.Cells.Item(94, 13).Value = -741.75
$ws.Cells.Item(107, 8).Value = 2826.7693
$ws.Cells.Item(107, 9).Value = 2625.7
$ws.Cells.Item(107, 10).Value = 3497
$ws.Cells.Item(107, 11).Value = 2625.7
$ws.Cells.Item(107, 12).Value = 3497
$ws.Cells.Item(107, 13).Value = -705.6999999999998
$ws.Cells.Item(107, 14).Value = -7337

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 1016.4667
$ws.Cells.Item(107, 10).Value = 1146.6
$ws.Cells.Item(107, 12).Value = 1146.6
$ws.Cells.Item(107, 14).Value = -4986.6
$ws.Cells.Item(122, 8).Value = 1409.8462
$ws.Cells.Item(122, 9).Value = 1501.2222
$ws.Cells.Item(122, 10).Value = 1204.25
$ws.Cells.Item(122, 11).Value = 4503.6666
$ws.Cells.Item(122, 12).Value = 3612.75
$ws.Cells.Item(122, 13).Value = -2053.6666
$ws.Cells.Item(122, 14).Value = -8512.75
$ws.Cells.Item(139, 8).Value = 58480
$ws.Cells.Item(139, 10).Value = 58480
$ws.Cells.Item(139, 12).Value = 58480
$ws.Cells.Item(139, 14).Value = -68760

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 2500799.5
$ws.Cells.Item(46, 9).Value = 800
$ws.Cells.Item(46, 10).Value = 3334132.8
$ws.Cells.Item(46, 11).Value = 2400
$ws.Cells.Item(46, 12).Value = 10002398.4
$ws.Cells.Item(46, 13).Value = -2309
$ws.Cells.Item(46, 14).Value = -10002580.4
$ws.Cells.Item(56, 8).Value = 10134.3
$ws.Cells.Item(56, 9).Value = 10134.3
$ws.Cells.Item(56, 11).Value = 10134.3
$ws.Cells.Item(56, 13).Value = -9604.299999999999
$ws.Cells.Item(107, 8).Value = 947.61536
$ws.Cells.Item(107, 10).Value = 1142.8334
$ws.Cells.Item(107, 12).Value = 3428.5002
$ws.Cells.Item(107, 14).Value = -7268.5002
$ws.Cells.Item(114, 8).Value = 27549.445
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = 27549.445
$ws.Cells.Item(114, 11).Value = 0
$ws.Cells.Item(114, 12).Value = 82648.33499999999
$ws.Cells.Item(114, 13).Value = $null
$ws.Cells.Item(114, 14).Value = -89156.33499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2165.7144
$ws.Cells.Item(102, 9).Value = 2224.6155
$ws.Cells.Item(102, 10).Value = 1400
$ws.Cells.Item(102, 11).Value = 2224.6155
$ws.Cells.Item(102, 12).Value = 1400
$ws.Cells.Item(102, 13).Value = -602.6154999999999
$ws.Cells.Item(102, 14).Value = -4644
$ws.Cells.Item(122, 8).Value = 51665.5
$ws.Cells.Item(122, 9).Value = 51665.5
$ws.Cells.Item(122, 11).Value = 154996.5
$ws.Cells.Item(122, 13).Value = -152546.5
$ws.Cells.Item(126, 8).Value = 11477067
$ws.Cells.Item(126, 9).Value = 6933.4287
$ws.Cells.Item(126, 10).Value = 19506160
$ws.Cells.Item(126, 11).Value = 20800.2861
$ws.Cells.Item(126, 12).Value = 58518480
$ws.Cells.Item(126, 13).Value = -18330.2861
$ws.Cells.Item(126, 14).Value = -58523420
$ws.Cells.Item(132, 8).Value = 8913.134
$ws.Cells.Item(132, 9).Value = 8252.166999999999
$ws.Cells.Item(132, 10).Value = 11557
$ws.Cells.Item(132, 11).Value = 24756.501
$ws.Cells.Item(132, 12).Value = 34671
$ws.Cells.Item(132, 13).Value = -22226.501
$ws.Cells.Item(132, 14).Value = -39731

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 45460070
$ws.Cells.Item(7, 9).Value = 5966.3335
$ws.Cells.Item(7, 10).Value = 100005000
$ws.Cells.Item(7, 11).Value = 5966.3335
$ws.Cells.Item(7, 12).Value = 100005000
$ws.Cells.Item(7, 13).Value = -5854.3335
$ws.Cells.Item(7, 14).Value = -100005224
$ws.Cells.Item(16, 8).Value = 987.2778
$ws.Cells.Item(16, 9).Value = 801.36365
$ws.Cells.Item(16, 10).Value = 1279.4286
$ws.Cells.Item(16, 11).Value = 801.36365
$ws.Cells.Item(16, 12).Value = 1279.4286
$ws.Cells.Item(16, 13).Value = -631.36365
$ws.Cells.Item(16, 14).Value = -1619.4286
$ws.Cells.Item(40, 8).Value = 4121.76
$ws.Cells.Item(40, 9).Value = 4310.864
$ws.Cells.Item(40, 10).Value = 2735
$ws.Cells.Item(40, 11).Value = 4310.864
$ws.Cells.Item(40, 12).Value = 2735
$ws.Cells.Item(40, 13).Value = -4174.864
$ws.Cells.Item(40, 14).Value = -3007
$ws.Cells.Item(126, 8).Value = 45460070
$ws.Cells.Item(126, 9).Value = 5966.3335
$ws.Cells.Item(126, 10).Value = 100005000
$ws.Cells.Item(126, 11).Value = 17899.0005
$ws.Cells.Item(126, 12).Value = 300015000
$ws.Cells.Item(126, 13).Value = -15429.0005
$ws.Cells.Item(126, 14).Value = -300019940
$ws.Cells.Item(132, 8).Value = 8666.121999999999
$ws.Cells.Item(132, 9).Value = 8356.393
$ws.Cells.Item(132, 10).Value = 9333.23
$ws.Cells.Item(132, 11).Value = 25069.179
$ws.Cells.Item(132, 12).Value = 27999.69
$ws.Cells.Item(132, 13).Value = -22539.179
$ws.Cells.Item(132, 14).Value = -33059.69
$ws.Cells.Item(135, 8).Value = 70000
$ws.Cells.Item(135, 10).Value = 70000
$ws.Cells.Item(135, 12).Value = 70000
$ws.Cells.Item(135, 14).Value = -80140

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 2803.3
$ws.Cells.Item(107, 9).Value = 1520
$ws.Cells.Item(107, 11).Value = 4560
$ws.Cells.Item(107, 13).Value = -2640
$ws.Cells.Item(113, 8).Value = 755.26666
$ws.Cells.Item(113, 10).Value = 522.5
$ws.Cells.Item(113, 12).Value = 1567.5
$ws.Cells.Item(113, 14).Value = -5907.5
$ws.Cells.Item(126, 8).Value = 95244024
$ws.Cells.Item(126, 10).Value = 166673660
$ws.Cells.Item(126, 12).Value = 500020980
$ws.Cells.Item(126, 14).Value = -500025920
